$d = $word.ActiveDocument

# Locate the paragraph that ends the "delete SEVIRI data" crontab line
# ("...del.log 2>/home/mariners/log/del-error.log") and then the blank
# paragraph that immediately follows it. The new crontab block (comment +
# command) is appended right after that blank paragraph, followed by a
# fresh blank paragraph of its own - matching the blank/comment/command
# rhythm used by every other entry earlier in the document - and all of
# it lands before the trailing, bookmark-only paragraph at the very end.
$findRange = $d.Content
$found = $findRange.Find.Execute("del-error.log", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$anchorPara = $findRange.Paragraphs(1)
$blankPara = $anchorPara.Next()

# Create three new empty paragraphs right after that blank paragraph:
# comment line, command line, trailing blank separator.
[void]$blankPara.Range.InsertParagraphAfter()
$commentPara = $blankPara.Next()
[void]$commentPara.Range.InsertParagraphAfter()
$commandPara = $commentPara.Next()
[void]$commandPara.Range.InsertParagraphAfter()
$trailingBlankPara = $commandPara.Next()

# Fill the comment/command paragraphs in with the exact WordprocessingML
# (runs + proofing marks) so the result matches what Word itself would
# have produced; the trailing separator paragraph is left empty.
$commentXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve"># </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>crontab</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> to download data from AOD MODIS-TERRA &amp; AQUA (10km)</w:t></w:r></w:p>'
[void]$commentPara.Range.InsertXML($commentXml)

$commandXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>30 17 * * * /home/mariners/MODIS_AOD/MODIS_AOD_download_hdf_ocean.sh 1&gt;/home/mariners/log/run_MODIS.log 2&gt;/home/mariners/log/run_MODIS-error.log</w:t></w:r></w:p>'
[void]$commandPara.Range.InsertXML($commandXml)

# Keep the trailing separator paragraph truly empty (no stray run), same
# as the blank paragraphs used elsewhere in the document.
$trailingBlankXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>'
[void]$trailingBlankPara.Range.InsertXML($trailingBlankXml)
